$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns.
# NumberFormat is forced to text ("@") before each write so that
# numeric-looking strings (e.g. "1.280", "1.003") are preserved exactly
# as text instead of being auto-coerced to numbers (which would drop
# significant trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.153.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5284'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.05%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2608'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06361'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07779'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.508'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.668.98'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5489'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅8210'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.56'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.156.11'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.003'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.589'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.051'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '141.42'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1251'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.287'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.24'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05944'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.280'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.528'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.583'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9549'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.791'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.414'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5718'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01619'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.818'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8487'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.07'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.026.57'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.802.06'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.39'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.009'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.496'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4296'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05155'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.841'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09724'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.21%  '
